# update scripts wuth new tpm
# Mif-Cd74 NATMI ligand-receptor pair sheet: refresh ligand/receptor expression
# values (and their derived specificity / edge-weight columns) with newly
# computed TPM figures for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.556445
$ws.Range("H2").Value = 19.669335
$ws.Range("I2").Value = 0.2003595613103873
$ws.Range("J2").Value = 0.2003595613103873
$ws.Range("M2").Value = 1.231278
$ws.Range("N2").Value = 3.693834
$ws.Range("O2").Value = 0.3283785416403858
$ws.Range("P2").Value = 0.3283785416403859
$ws.Range("Q2").Value = 8.072806486709998
$ws.Range("R2").Value = 72.65525838038999
$ws.Range("S2").Value = 0.06579378054681244
$ws.Range("T2").Value = 0.06579378054681245
$ws.Range("G3").Value = 6.556445
$ws.Range("H3").Value = 19.669335
$ws.Range("I3").Value = 0.2003595613103873
$ws.Range("J3").Value = 0.2003595613103873
$ws.Range("O3").Value = 0.2084514246837437
$ws.Range("P3").Value = 0.2084514246837437
$ws.Range("Q3").Value = 5.124537081335
$ws.Range("R3").Value = 46.120833732015
$ws.Range("S3").Value = 0.04176523600416011
$ws.Range("T3").Value = 0.04176523600416011
$ws.Range("G4").Value = 6.556445
$ws.Range("H4").Value = 19.669335
$ws.Range("I4").Value = 0.2003595613103873
$ws.Range("J4").Value = 0.2003595613103873
$ws.Range("O4").Value = 0.4631700336758705
$ws.Range("P4").Value = 0.4631700336758705
$ws.Range("Q4").Value = 11.38649935416
$ws.Range("R4").Value = 102.47849418744
$ws.Range("S4").Value = 0.0928005447594147
$ws.Range("T4").Value = 0.09280054475941471
$ws.Range("I5").Value = 0.3842320902647997
$ws.Range("J5").Value = 0.3842320902647997
$ws.Range("M5").Value = 1.231278
$ws.Range("N5").Value = 3.693834
$ws.Range("O5").Value = 0.3283785416403858
$ws.Range("P5").Value = 0.3283785416403859
$ws.Range("Q5").Value = 15.48132412751
$ws.Range("R5").Value = 139.33191714759
$ws.Range("S5").Value = 0.126173573452592
$ws.Range("T5").Value = 0.126173573452592
$ws.Range("I6").Value = 0.3842320902647997
$ws.Range("J6").Value = 0.3842320902647997
$ws.Range("O6").Value = 0.2084514246837437
$ws.Range("P6").Value = 0.2084514246837437
$ws.Range("S6").Value = 0.08009372662491028
$ws.Range("T6").Value = 0.08009372662491029
$ws.Range("I7").Value = 0.3842320902647997
$ws.Range("J7").Value = 0.3842320902647997
$ws.Range("O7").Value = 0.4631700336758705
$ws.Range("P7").Value = 0.4631700336758705
$ws.Range("S7").Value = 0.1779647901872974
$ws.Range("T7").Value = 0.1779647901872974
$ws.Range("G8").Value = 13.59357133333334
$ws.Range("I8").Value = 0.4154083484248129
$ws.Range("J8").Value = 0.415408348424813
$ws.Range("M8").Value = 1.231278
$ws.Range("N8").Value = 3.693834
$ws.Range("O8").Value = 0.3283785416403858
$ws.Range("P8").Value = 0.3283785416403859
$ws.Range("Q8").Value = 16.737465324164
$ws.Range("R8").Value = 150.637187917476
$ws.Range("S8").Value = 0.1364111876409813
$ws.Range("T8").Value = 0.1364111876409814
$ws.Range("G9").Value = 13.59357133333334
$ws.Range("I9").Value = 0.4154083484248129
$ws.Range("J9").Value = 0.415408348424813
$ws.Range("O9").Value = 0.2084514246837437
$ws.Range("P9").Value = 0.2084514246837437
$ws.Range("Q9").Value = 10.62477613484734
$ws.Range("R9").Value = 95.62298521362601
$ws.Range("S9").Value = 0.08659246205467325
$ws.Range("T9").Value = 0.08659246205467326
$ws.Range("G10").Value = 13.59357133333334
$ws.Range("I10").Value = 0.4154083484248129
$ws.Range("J10").Value = 0.415408348424813
$ws.Range("O10").Value = 0.4631700336758705
$ws.Range("P10").Value = 0.4631700336758705
$ws.Range("S10").Value = 0.1924046987291583
$ws.Range("T10").Value = 0.1924046987291584
